$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 4 data, reusing the same number formats as rows 2/3
$ws.Range("A4").Value = 44317
$ws.Range("A4").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("B4").Value = 0.16666666666666666
$ws.Range("B4").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("C4").Value = "apparence - modification affihchage JS"

# Update selection to C5 like in the target file
$ws.Range("C5").Select()
